$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calculator")

# Update the formula in AD7 with the re-fit coefficients
$ws.Range("AD7").Formula = "=10^(1.071037240663 + 2.444855540303*Z7 + 0.090107006387*T7 + -0.895938092707*Z7^2 + -0.182764652802*Z7*T7 + -0.002205819859*T7^2 + 0.117839927217*Z7^2*T7 + 0.005369514927*Z7*T7^2 + -0.003515592778*Z7^2*T7^2)"

# Add a new row describing the fitted equation, used as a documentation / check row
$ws.Range("T9").Value = "log10(pCO2) = "
$ws.Range("U9").Value = "(1.071037240663 +  2.444855540303*Z7 +  0.090107006387*T7 +  -0.895938092707*Z7^2 +  -0.182764652802*Z7*T7 +  -0.002205819859*T7^2 +  0.117839927217*Z7^2*T7 +  0.005369514927*Z7*T7^2 +  -0.003515592778*Z7^2*T7^2)"

# Widen the input/output columns to comfortably fit the new content
$ws.Columns.Item(20).ColumnWidth = 14.8
$ws.Columns.Item(26).ColumnWidth = 14.8
$ws.Columns.Item(30).ColumnWidth = 14.8
